# Data structures module is completed. changes in home page too.
#
# 1) Add a new "DataStructures" worksheet at the end of the workbook and
#    populate it with the python-code / result table.
# 2) Make the new sheet the active sheet (bookView activeTab shifts to it,
#    and the previously-active "linkedList" sheet loses tabSelected).

$wb = $excel.ActiveWorkbook

$signindata  = $wb.Worksheets.Item(1)
$linkedList  = $wb.Worksheets.Item(2)

# Add the new sheet after the last existing sheet so it lands at the end.
$ws = $wb.Worksheets.Add($null, $linkedList)
$ws.Name = "DataStructures"

# ---- Row 1 : header ----------------------------------------------------
$ws.Range("A1").Value = "pythoncode"
$ws.Range("A1").Font.Name = "Calibri"
$ws.Range("A1").Font.Size = 20
$ws.Range("A1").Font.Color = 0

$ws.Range("B1").Value = "Result"
$ws.Range("B1").Font.Name = "Calibri"
$ws.Range("B1").Font.Size = 16
$ws.Range("B1").Font.Color = 0

# ---- Row 2 : print("hello");  ->  hello --------------------------------
$ws.Range("A2").Value = 'print("hello");'
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.Size = 16
$ws.Range("A2").Font.Color = 0

$ws.Range("B2").Value = "hello"
$ws.Range("B2").Font.Name = "Calibri"
$ws.Range("B2").Font.Size = 16
$ws.Range("B2").Font.Color = 0

# ---- Row 3 : print("welcome)";  ->  SyntaxError ------------------------
$ws.Range("A3").Value = 'print("welcome)”;'
$ws.Range("A3").Font.Name = "Calibri"
$ws.Range("A3").Font.Size = 16
$ws.Range("A3").Font.Color = 0

$ws.Range("B3").Value = "SyntaxError: bad input on line 1"
$ws.Range("B3").Font.Name = "Calibri"
$ws.Range("B3").Font.Size = 14
$ws.Range("B3").Font.Color = 0

# ---- Row 4 : def search(...)  ->  Element Found ------------------------
$ws.Range("A4").Value = "def search(lst,value):" + [char]10 + "if(value in lst):" + [char]10 + 'return "Element Found"' + [char]10 + "else:" + [char]10 + 'return "Not Found"' + [char]10 + "pass" + [char]10 + "print(search([1,2,3],1))"

# Copy formatting (fill + left/wrap alignment) from the existing highlighted
# cell on "linkedList" so the new style reuses the same fill, then bump the
# font to match.
$linkedList.Range("A4").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Font.Name = "Calibri"
$ws.Range("A4").Font.Size = 14
$ws.Range("A4").Font.Color = 0
$ws.Range("A4").WrapText = $true

$ws.Range("B4").Value = "Element Found"
$ws.Range("B4").Font.Name = "Calibri"
$ws.Range("B4").Font.Size = 16
$ws.Range("B4").Font.Color = 0

# ---- Row heights / column widths ---------------------------------------
$ws.Rows.Item(1).RowHeight = 24.45
$ws.Rows.Item(2).RowHeight = 19.7
$ws.Rows.Item(3).RowHeight = 19.7
$ws.Rows.Item(4).RowHeight = 119.4

$ws.Columns.Item(1).ColumnWidth = 27.15
$ws.Columns.Item(2).ColumnWidth = 27.29

# ---- Page setup / margins / header-footer ------------------------------
$ws.PageSetup.LeftMargin = 56.7
$ws.PageSetup.RightMargin = 56.7
$ws.PageSetup.TopMargin = 75.8
$ws.PageSetup.BottomMargin = 75.8
$ws.PageSetup.HeaderMargin = 56.7
$ws.PageSetup.FooterMargin = 56.7
$ws.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&Kffffff&A'
$ws.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12&KffffffPage &P'

# ---- Selection / view ----------------------------------------------------
# Activating the new sheet makes it the active tab; "linkedList" (previously
# tabSelected) automatically loses that flag.
$ws.Activate() | Out-Null
$ws.Range("C4").Select() | Out-Null
$excel.ActiveWindow.Zoom = 80
